# Pontuacoes sheet was reshaped:
#   - old layout: B1="Nome", C1="Pontuação", D1="Unnamed: 0" headers (row 1,
#     styled) with two data rows (A/B/C/D, row 2-3) plus an index column.
#   - new layout: A1="Nome", B1="Pontos" headers (styled) with a single data
#     row A2="Diego", B2=40.
#
# NOTE: this headless COM host does not keep the clipboard alive across a
# Range.Clear() call (even on an unrelated range), so Copy()/PasteSpecial()
# cannot be used to carry the header style through the clear. Instead we
# rebuild the bold/bordered/centered header look with direct formatting
# calls after writing the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old A1:D3 content (values + styles) entirely.
$ws.Range("A1:D3").Clear()

# New header row.
$ws.Cells.Item(1, 1).Value = "Nome"
$ws.Cells.Item(1, 2).Value = "Pontos"

# Re-apply the header formatting (bold, thin box border, centered/top aligned)
# that the original "Nome"/"Pontuação"/"Unnamed: 0" header cells carried.
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop

# New (single) data row, left unstyled like the original data cells.
$ws.Cells.Item(2, 1).Value = "Diego"
$ws.Cells.Item(2, 2).Value = 40
